$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 12502343
$ws.Range("J17").Value = 14288249
$ws.Range("L17").Value = 42864747
$ws.Range("N17").Value = -42865083
$ws.Range("H32").Value = 4227.2856
$ws.Range("I32").Value = 4098
$ws.Range("K32").Value = 4098
$ws.Range("M32").Value = -3772
$ws.Range("H42").Value = 1000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H69").Value = 5250
$ws.Range("I69").Value = 3000
$ws.Range("K69").Value = 9000
$ws.Range("M69").Value = -8126
$ws.Range("H72").Value = 5250
$ws.Range("I72").Value = 3000
$ws.Range("K72").Value = 27000
$ws.Range("M72").Value = -22632
$ws.Range("H80").Value = 359.2857
$ws.Range("I80").Value = 352.63635
$ws.Range("K80").Value = 1057.90905
$ws.Range("M80").Value = -59.90904999999998
$ws.Range("H83").Value = 359.2857
$ws.Range("I83").Value = 352.63635
$ws.Range("K83").Value = 3173.72715
$ws.Range("M83").Value = 1818.27285
$ws.Range("H103").Value = 537.25
$ws.Range("I103").Value = 499
$ws.Range("J103").Value = 550
$ws.Range("K103").Value = 1497
$ws.Range("L103").Value = 1650
$ws.Range("M103").Value = -911
$ws.Range("N103").Value = -2822
$ws.Range("H127").Value = 1097.4
$ws.Range("I127").Value = 1121.75
$ws.Range("K127").Value = 3365.25
$ws.Range("M127").Value = 1594.75
$ws.Range("H133").Value = 77322.25
$ws.Range("J133").Value = 77322.25
$ws.Range("L133").Value = 77322.25
$ws.Range("N133").Value = -87442.25
$ws.Range("H138").Value = 1969.1818
$ws.Range("I138").Value = 1275.5238
$ws.Range("J138").Value = 2397.6177
$ws.Range("K138").Value = 3826.5714
$ws.Range("L138").Value = 7192.853099999999
$ws.Range("M138").Value = 1313.4286
$ws.Range("N138").Value = -17472.8531
$ws.Range("H141").Value = 4853.2
$ws.Range("I141").Value = 4853.2
$ws.Range("K141").Value = 14559.6
$ws.Range("M141").Value = -9379.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1497.8125
$ws.Range("I32").Value = 1523.6129
$ws.Range("K32").Value = 1523.6129
$ws.Range("M32").Value = -1236.6129
$ws.Range("H102").Value = 4824.9546
$ws.Range("I102").Value = 3619.6667
$ws.Range("J102").Value = 10248.75
$ws.Range("K102").Value = 3619.6667
$ws.Range("L102").Value = 10248.75
$ws.Range("M102").Value = -1997.6667
$ws.Range("N102").Value = -13492.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1749.3928
$ws.Range("I94").Value = 1410.2778
$ws.Range("J94").Value = 2359.8
$ws.Range("K94").Value = 1410.2778
$ws.Range("L94").Value = 2359.8
$ws.Range("M94").Value = -959.2778000000001
$ws.Range("N94").Value = -3261.8
$ws.Range("H99").Value = 4990.727
$ws.Range("I99").Value = 3664
$ws.Range("K99").Value = 3664
$ws.Range("M99").Value = -2166
$ws.Range("H107").Value = 1743.5
$ws.Range("I107").Value = 1654.4445
$ws.Range("K107").Value = 1654.4445
$ws.Range("M107").Value = 265.5554999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9736.066000000001
$ws.Range("I58").Value = 6181.6665
$ws.Range("J58").Value = 12105.667
$ws.Range("K58").Value = 6181.6665
$ws.Range("L58").Value = 12105.667
$ws.Range("M58").Value = -5978.6665
$ws.Range("N58").Value = -12511.667
$ws.Range("H68").Value = 41776.332
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 41776.332
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H99").Value = 3608.5334
$ws.Range("I99").Value = 3943.111
$ws.Range("J99").Value = 3106.6667
$ws.Range("K99").Value = 3943.111
$ws.Range("L99").Value = 3106.6667
$ws.Range("M99").Value = -2445.111
$ws.Range("N99").Value = -6102.6667
$ws.Range("H122").Value = 3664.2273
$ws.Range("I122").Value = 3436.8
$ws.Range("J122").Value = 4151.5713
$ws.Range("K122").Value = 10310.4
$ws.Range("L122").Value = 12454.7139
$ws.Range("M122").Value = -7860.400000000001
$ws.Range("N122").Value = -17354.7139
$ws.Range("H126").Value = 3608.5334
$ws.Range("I126").Value = 3943.111
$ws.Range("J126").Value = 3106.6667
$ws.Range("K126").Value = 11829.333
$ws.Range("L126").Value = 9320.000100000001
$ws.Range("M126").Value = -9359.332999999999
$ws.Range("N126").Value = -14260.0001
$ws.Range("H136").Value = 9736.066000000001
$ws.Range("I136").Value = 6181.6665
$ws.Range("J136").Value = 12105.667
$ws.Range("K136").Value = 18544.9995
$ws.Range("L136").Value = 36317.001
$ws.Range("M136").Value = -15994.9995
$ws.Range("N136").Value = -41417.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1503.8572
$ws.Range("J18").Value = 1749.5
$ws.Range("L18").Value = 5248.5
$ws.Range("N18").Value = -5586.5
$ws.Range("H33").Value = 565.625
$ws.Range("I33").Value = 287.8
$ws.Range("J33").Value = 1028.6666
$ws.Range("K33").Value = 1726.8
$ws.Range("L33").Value = 6171.9996
$ws.Range("M33").Value = -1443.8
$ws.Range("N33").Value = -6737.9996
$ws.Range("H107").Value = 1325.1666
$ws.Range("J107").Value = 1533.6
$ws.Range("L107").Value = 4600.799999999999
$ws.Range("N107").Value = -8440.799999999999
$ws.Range("H131").Value = 22730826
$ws.Range("J131").Value = 3957.5789
$ws.Range("L131").Value = 11872.7367
$ws.Range("N131").Value = -21952.7367
$ws.Range("H137").Value = 2137.6
$ws.Range("I137").Value = 1735.6923
$ws.Range("K137").Value = 5207.0769
$ws.Range("M137").Value = -107.0769
$ws.Range("H139").Value = 2093.7778
$ws.Range("I139").Value = 1865.9333
$ws.Range("K139").Value = 5597.7999
$ws.Range("M139").Value = -457.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 759000
$ws.Range("J21").Value = 18000
$ws.Range("L21").Value = 18000
$ws.Range("N21").Value = -18346
$ws.Range("H24").Value = 15338
$ws.Range("J24").Value = 16005.25
$ws.Range("L24").Value = 16005.25
$ws.Range("N24").Value = -16351.25
$ws.Range("H30").Value = 759000
$ws.Range("J30").Value = 18000
$ws.Range("L30").Value = 18000
$ws.Range("N30").Value = -18210
$ws.Range("H97").Value = 1168.6666
$ws.Range("I97").Value = 853.8333
$ws.Range("J97").Value = 1798.3334
$ws.Range("K97").Value = 853.8333
$ws.Range("L97").Value = 1798.3334
$ws.Range("M97").Value = -357.8333
$ws.Range("N97").Value = -2790.3334
$ws.Range("H132").Value = 2394.6875
$ws.Range("I132").Value = 2354.3333
$ws.Range("K132").Value = 7062.999899999999
$ws.Range("M132").Value = -4532.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3905.875
$ws.Range("J7").Value = 3999.5
$ws.Range("L7").Value = 3999.5
$ws.Range("N7").Value = -4223.5
$ws.Range("H82").Value = 699.2
$ws.Range("I82").Value = 699.3333
$ws.Range("K82").Value = 699.3333
$ws.Range("M82").Value = -338.3333
$ws.Range("H85").Value = 699.2
$ws.Range("I85").Value = 699.3333
$ws.Range("K85").Value = 699.3333
$ws.Range("M85").Value = 548.6667
$ws.Range("H126").Value = 3905.875
$ws.Range("J126").Value = 3999.5
$ws.Range("L126").Value = 11998.5
$ws.Range("N126").Value = -16938.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9496
$ws.Range("J15").Value = 9495
$ws.Range("L15").Value = 9495
$ws.Range("N15").Value = -10071
$ws.Range("H81").Value = 1924.8572
$ws.Range("J81").Value = 4641
$ws.Range("L81").Value = 9282
$ws.Range("N81").Value = -11404
$ws.Range("H84").Value = 1924.8572
$ws.Range("J84").Value = 4641
$ws.Range("L84").Value = 46410
$ws.Range("N84").Value = -57018
$ws.Range("H132").Value = 7957.75
$ws.Range("I132").Value = 7386.625
$ws.Range("J132").Value = 9100
$ws.Range("K132").Value = 22159.875
$ws.Range("L132").Value = 27300
$ws.Range("M132").Value = -19629.875
$ws.Range("N132").Value = -32360
